$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 currently holds the shared string "R40". It needs to become the
# text "1" (a new shared-string entry), while keeping the cell's existing
# style/format untouched (no number-format or style-table changes).
#
# Assigning Range.Value directly with a numeric-looking string ("1")
# would make Excel auto-convert it to a number, losing the text type.
# To force a genuine text value without touching NumberFormat (which
# would allocate a brand-new cell style), write it as a text formula
# result and then paste just the value back over itself.
$cell = $ws.Range("B11")
$cell.Formula = '="1"'
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
$excel.CutCopyMode = 0
